# Changed to ContractUptime test-case
# Add a new "DownTime" column (W) with a numeric value of 360 in row 2,
# matching the svmx an_Datasheet1.xlsx test-data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (W1) - shared string "DownTime"
$ws.Range("W1").Value = "DownTime"

# New data cell (W2) - numeric value 360
$ws.Range("W2").Value = 360

# Give the new column a sensible custom width, matching the other
# data columns in this sheet.
$ws.Columns.Item(23).ColumnWidth = 27.5

# Keep the selection on the newly added cell, same as the original edit.
[void]$ws.Range("W2").Select()
